$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) semantics for numeric-looking values in columns D and E
# by temporarily applying a text number format, then restore the default
# "Normal" style afterwards so the saved cells carry no style index,
# matching how the sheet stores its other plain text cells.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range('D2').Value = '306.38'
$ws.Range('E2').Value = '-0.66%'

$ws.Range('D3').Value = '38.90'
$ws.Range('E3').Value = '7.19%'

$ws.Range('D4').Value = '5.114'
$ws.Range('E4').Value = '1.32%'

$ws.Range('D5').Value = '0.08085'
$ws.Range('E5').Value = '-0.67%'

$ws.Range('D6').Value = '1.934'
$ws.Range('E6').Value = '-6.40%'

$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = '4.184'
$ws.Range('E7').Value = '0.63%'

$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D8').Value = '7.981'
$ws.Range('E8').Value = '1.34%'

$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '0.9297'
$ws.Range('E9').Value = '0.12%'

$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '0.1473'
$ws.Range('E10').Value = '-0.10%'

$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1927'
$ws.Range('E11').Value = '-0.29%'

$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.09040'
$ws.Range('E12').Value = '-0.75%'

$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.03504'
$ws.Range('E13').Value = '1.36%'

$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.09781'
$ws.Range('E14').Value = '-1.17%'

$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001403'
$ws.Range('E15').Value = '-0.56%'

$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.006066'
$ws.Range('E16').Value = '-5.61%'

$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.782'
$ws.Range('E17').Value = '-1.59%'

$ws.Range('D18').Value = '3.451'
$ws.Range('E18').Value = '1.49%'

$ws.Range('E19').Value = '-0.18%'

$ws.Range('D20').Value = '0.1348'
$ws.Range('E20').Value = '2.17%'

$ws.Range('D21').Value = '4.683'
$ws.Range('E21').Value = '-2.65%'

$ws.Range('D22').Value = '0.2415'
$ws.Range('E22').Value = '3.36%'

$ws.Range('D23').Value = '0.04380'
$ws.Range('E23').Value = '-0.04%'

$ws.Range('E24').Value = '0.34%'

$ws.Range('D25').Value = '0.004283'
$ws.Range('E25').Value = '2.15%'

$ws.Range('D26').Value = '0.0001303'
$ws.Range('E26').Value = '0.36%'

$ws.Range('D39').Value = '0.02041'
$ws.Range('E39').Value = '-0.11%'

$ws.Range('E40').Value = '-2.10%'

$ws.Range('D41').Value = '0.007623'
$ws.Range('E41').Value = '2.11%'

$ws.Range('D42').Value = '0.01011'
$ws.Range('E42').Value = '-0.34%'

$ws.Range('D43').Value = '0.1350'
$ws.Range('E43').Value = '-2.06%'

$ws.Range('D44').Value = '0.002125'
$ws.Range('E44').Value = '-0.10%'

$ws.Range('D45').Value = '0.009916'
$ws.Range('E45').Value = '2.48%'

$ws.Range('D46').Value = '0.00006185'
$ws.Range('E46').Value = '-2.05%'

$ws.Range('E47').Value = '0.24%'

$ws.Range('D48').Value = '0.002891'

$ws.Range('E49').Value = '0.13%'

$ws.Range('E50').Value = '0.24%'

$ws.Range('E51').Value = '0.24%'

# Restore the default style on the numeric-text columns so no stray
# quote-prefix / text-format style sticks to the cells.
$numRange.Style = "Normal"
